$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-03-26 Wednesday" "2025-03-27 Thursday"

Replace-Text "486×7=" "499×3="
Replace-Text "520×5=" "613×3="
Replace-Text "667×6=" "141×8="
Replace-Text "953×6=" "977×8="
Replace-Text "225×5=" "221×6="
Replace-Text "965×9=" "198×8="
Replace-Text "920×6=" "806×5="
Replace-Text "684×6=" "610×7="
Replace-Text "310×9=" "819×8="
Replace-Text "972×7=" "517×8="
Replace-Text "852×8=" "391×2="
Replace-Text "386×4=" "881×6="
Replace-Text "245×7=" "538×3="
Replace-Text "894×6=" "692×4="
Replace-Text "139×4=" "314×6="
Replace-Text "623×9=" "473×9="
Replace-Text "121×5=" "152×5="
Replace-Text "177×8=" "213×2="
Replace-Text "690×9=" "974×6="
Replace-Text "373×8=" "829×3="
Replace-Text "975×2=" "954×2="
Replace-Text "810×6=" "430×5="
Replace-Text "345×6=" "607×9="
Replace-Text "864×8=" "698×5="
Replace-Text "557×5=" "583×5="
